$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Min_Money" column (E) entirely; this shifts Peek/Min_Peek/Trade/Min_Trade left.
$ws.Range("E1").EntireColumn.Delete()

# Rename the headers that moved into E/F/H to reflect the new "Ally Number" game data.
$ws.Range("E1").Value = "Ally Number"
$ws.Range("F1").Value = "Min Ally Number"
$ws.Range("H1").Value = "Min Trade"

# Match the workbook's recorded active selection after the edit.
$ws.Range("H2").Select() | Out-Null
